$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "41.767.07"
Set-TextValue "E2" "  -1.78%  "

Set-TextValue "D3" "2.216.06"
Set-TextValue "E3" "  -1.71%  "

Set-TextValue "E4" "  +0.28%  "

Set-TextValue "D5" "241.22"
Set-TextValue "E5" "  -2.16%  "

Set-TextValue "E6" "  -1.15%  "

Set-TextValue "D7" "72.43"
Set-TextValue "E7" "  -5.65%  "

Set-TextValue "E8" "  +0.12%  "

Set-TextValue "E9" "  -4.59%  "

Set-TextValue "D10" "41.98"
Set-TextValue "E10" "  -7.03%  "

Set-TextValue "D11" "0.0948"
Set-TextValue "E11" "  -0.40%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D12" "6.95"
Set-TextValue "E12" "  -4.51%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D13" "0.103"
Set-TextValue "E13" "  +0.30%  "

Set-TextValue "D14" "2.549.21"
Set-TextValue "E14" "  -1.72%  "

Set-TextValue "D15" "14.24"
Set-TextValue "E15" "  -3.02%  "

Set-TextValue "E16" "  -3.19%  "

Set-TextValue "D17" "2.217.11"
Set-TextValue "E17" "  -3.16%  "

Set-TextValue "D18" "41.637.46"
Set-TextValue "E18" "  -1.70%  "

Set-TextValue "E19" "  +3.52%  "

Set-TextValue "D20" "72.45"
Set-TextValue "E20" "  +0.28%  "

Set-TextValue "E21" "  -0.99%  "

Set-TextValue "D22" "11.06"
Set-TextValue "E22" "  +22.12%  "

Set-TextValue "D23" "229.32"
Set-TextValue "E23" "  -1.19%  "

Set-TextValue "E24" "  -9.88%  "

Set-TextValue "E25" "  +0.02%  "

Set-TextValue "E26" "  -1.69%  "

Set-TextValue "E27" "  +0.30%  "

Set-TextValue "E28" "  -1.83%  "

Set-TextValue "E29" "  -0.75%  "

Set-TextValue "D30" "167.30"
Set-TextValue "E30" "  -0.16%  "

Set-TextValue "D31" "20.43"
Set-TextValue "E31" "  -1.40%  "

Set-TextValue "D32" "0.0796"
Set-TextValue "E32" "  -3.66%  "

Set-TextValue "D33" "5.48"
Set-TextValue "E33" "  +3.29%  "

Set-TextValue "D34" "29.96"
Set-TextValue "E34" "  -3.90%  "

Set-TextValue "E35" "  -0.68%  "

Set-TextValue "D36" "0.107"
Set-TextValue "E36" "  -10.84%  "

Set-TextValue "D37" "4.28"
Set-TextValue "E37" "  -6.65%  "

Set-TextValue "E38" "  -4.52%  "

Set-TextValue "D39" "13.28"
Set-TextValue "E39" "  -6.41%  "

Set-TextValue "D40" "2.12"
Set-TextValue "E40" "  -3.15%  "

$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue "D41" "64.03"
Set-TextValue "E41" "  +0.17%  "

$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue "D42" "5.60"
Set-TextValue "E42" "  -3.96%  "

Set-TextValue "E43" "  -2.56%  "

Set-TextValue "D44" "8.71"
Set-TextValue "E44" "  -1.10%  "

Set-TextValue "D45" "102.73"
Set-TextValue "E45" "  -4.92%  "

Set-TextValue "D46" "0.0999"
Set-TextValue "E46" "  -3.24%  "

Set-TextValue "D47" "2.32"
Set-TextValue "E47" "  -2.71%  "

$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D48" "1.16"
Set-TextValue "E48" "  -2.38%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D49" "1.10"
Set-TextValue "E49" "  -3.00%  "

Set-TextValue "E50" "  -0.71%  "

Set-TextValue "D51" "2.423.61"
Set-TextValue "E51" "  -1.67%  "
